$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D (Price) holds locale-formatted text values (e.g. thousand-dot
# separators) in the source data, so force text formatting before writing
# so Excel does not reinterpret/renormalize them as numbers.
$ws.Range('D2:D51').NumberFormat = '@'

$ws.Range('D2').Value = '69.265.73'
$ws.Range('E2').Value = '  +0.19%  '
$ws.Range('D3').Value = '3.867.90'
$ws.Range('E3').Value = '  +3.24%  '
$ws.Range('E4').Value = '  -0.08%  '
$ws.Range('D5').Value = '604.01'
$ws.Range('E5').Value = '  +0.25%  '
$ws.Range('D6').Value = '164.92'
$ws.Range('E6').Value = '  -2.14%  '
$ws.Range('D7').Value = '3.867.60'
$ws.Range('E7').Value = '  +3.31%  '
$ws.Range('E8').Value = '  -0.03%  '
$ws.Range('E9').Value = '  -1.35%  '
$ws.Range('E10').Value = '  +0.51%  '
$ws.Range('E11').Value = '  -1.73%  '
$ws.Range('E12').Value = '  +0.24%  '
$ws.Range('D13').Value = '37.25'
$ws.Range('E13').Value = '  -2.51%  '
$ws.Range('E14').Value = '  -0.97%  '
$ws.Range('D15').Value = '4.500.54'
$ws.Range('E15').Value = '  +2.86%  '
$ws.Range('D16').Value = '3.837.91'
$ws.Range('E16').Value = '  +2.49%  '
$ws.Range('D17').Value = '69.359.83'
$ws.Range('E17').Value = '  +0.30%  '
$ws.Range('D18').Value = '7.63'
$ws.Range('E18').Value = '  +3.95%  '
$ws.Range('D19').Value = '11.76'
$ws.Range('E19').Value = '  +6.44%  '
$ws.Range('E20').Value = '  +0.23%  '
$ws.Range('D21').Value = '17.24'
$ws.Range('E21').Value = '  +0.50%  '
$ws.Range('D22').Value = '490.28'
$ws.Range('E22').Value = '  -0.83%  '
$ws.Range('D23').Value = '0.725'
$ws.Range('E23').Value = '  +0.19%  '
$ws.Range('E24').Value = '  +5.29%  '
$ws.Range('D25').Value = '84.69'
$ws.Range('D26').Value = '2.28'
$ws.Range('E26').Value = '  -1.69%  '
$ws.Range('D27').Value = '12.27'
$ws.Range('E27').Value = '  -0.68%  '
$ws.Range('E28').Value = '  -1.38%  '
$ws.Range('E29').Value = '  +0.12%  '
$ws.Range('D30').Value = '2.99'
$ws.Range('E30').Value = '  -0.16%  '
$ws.Range('D31').Value = '7.99'
$ws.Range('E31').Value = '  -0.55%  '
$ws.Range('B32').Value = 'EthereumClassic'
$ws.Range('C32').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D32').Value = '32.65'
$ws.Range('E32').Value = '  +3.24%  '
$ws.Range('B33').Value = 'ImmutableX'
$ws.Range('C33').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D33').Value = '2.40'
$ws.Range('E33').Value = '  -3.50%  '
$ws.Range('B34').Value = 'WrappedeETH'
$ws.Range('C34').Value = 'https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth'
$ws.Range('D34').Value = '4.018.43'
$ws.Range('E34').Value = '  +3.22%  '
$ws.Range('D35').Value = '3.809.08'
$ws.Range('E35').Value = '  +3.42%  '
$ws.Range('E36').Value = '  -1.04%  '
$ws.Range('D37').Value = '1.04'
$ws.Range('E37').Value = '  +1.83%  '
$ws.Range('E38').Value = '  +4.29%  '
$ws.Range('D39').Value = '5.94'
$ws.Range('E39').Value = '  +1.03%  '
$ws.Range('D40').Value = '0.998'
$ws.Range('E40').Value = '  -0.22%  '
$ws.Range('D41').Value = '0.321'
$ws.Range('E41').Value = '  -1.12%  '
$ws.Range('D42').Value = '3.03'
$ws.Range('E42').Value = '  +1.62%  '
$ws.Range('D43').Value = '438.69'
$ws.Range('E43').Value = '  +1.19%  '
$ws.Range('B44').Value = 'Stacks'
$ws.Range('C44').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D44').Value = '2.00'
$ws.Range('E44').Value = '  +0.02%  '
$ws.Range('B45').Value = 'OKB'
$ws.Range('C45').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D45').Value = '48.60'
$ws.Range('E45').Value = '  +0.18%  '
$ws.Range('D46').Value = '8.45'
$ws.Range('E46').Value = '  -0.80%  '
$ws.Range('E47').Value = '  +0.00%  '
$ws.Range('D48').Value = '27.72'
$ws.Range('E48').Value = '  +19.10%  '
$ws.Range('B49').Value = 'Maker'
$ws.Range('C49').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D49').Value = '2.857.16'
$ws.Range('E49').Value = '  +2.28%  '
$ws.Range('B50').Value = 'Monero'
$ws.Range('C50').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D50').Value = '143.09'
$ws.Range('E50').Value = '  +1.26%  '
$ws.Range('D51').Value = '0.0358'
$ws.Range('E51').Value = '  +1.47%  '
